# feat: add 2022-Q1 data
#
# Before: sheets are 2021-Q1, 2021-Q2, 2021-Q3, 2021-Q4, 总计 (totals).
# After:  a new "2022-Q1" fund-holdings detail sheet is inserted where the
#         old "总计" sheet used to be (keeping its sheetId/rId), and a
#         fresh "总计" sheet is appended right after it, carrying the same
#         totals table as before plus a new leading 2022-Q1 row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Step 1: turn the existing "总计" sheet into the new "2022-Q1" detail sheet
# ---------------------------------------------------------------------------
$q1 = $wb.Worksheets.Item(5)
$q1.Name = "2022-Q1"
$q1.Cells.ClearContents()

# Extend the bold/bordered header style (already present on B1:D1) across
# to H1 so every header cell shares the same formatting.
$q1.Cells.Item(1, 2).Copy()
$q1.Range("E1:H1").PasteSpecial(-4122)

# Extend the index-column style (already present on A2:A5) down to A13.
$q1.Cells.Item(2, 1).Copy()
$q1.Range("A6:A13").PasteSpecial(-4122)

$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($col = 0; $col -lt $headers.Length; $col++) {
    $q1.Cells.Item(1, $col + 2).Value = $headers[$col]
}

# code, name, scale, stock-position-total, position-ratio, market-value, rank
$rows = @(
    ,@("001128","宝盈新兴产业灵活配置混合","19.94","86.76","4.05","0.8076",6)
    ,@("320006","诺安灵活配置混合","10.15","70.73","7.33","0.7440",2)
    ,@("001702","东方创新科技混合","7.27","92.90","7.41","0.5387",2)
    ,@("001487","宝盈优势产业灵活配置混合","17.02","91.61","3.11","0.5293",7)
    ,@("001877","宝盈国家安全战略沪港深股票","13.59","90.66","3.67","0.4988",8)
    ,@("160642","鹏华增瑞灵活配置混合(LOF)","6.76","91.34","5.35","0.3617",5)
    ,@("959991","兴证资管金麒麟领先优势一年持有期混合A","8.28","79.70","3.96","0.3279",9)
    ,@("160919","大成产业升级股票(LOF)","3.95","87.76","5.79","0.2287",2)
    ,@("008988","大成科技创新混合A","2.65","91.66","5.59","0.1481",6)
    ,@("008989","大成科技创新混合C","1.11","91.66","5.59","0.0620",6)
    ,@("090009","大成行业轮动混合","1.71","82.88","2.33","0.0398",10)
    ,@("959993","兴证资管金麒麟领先优势一年持有期混合C","0.43","79.70","3.96","0.0170",9)
)

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $i + 2
    $row = $rows[$i]
    $q1.Cells.Item($r, 1).Value = $i
    # Fund code / scale / position-total / position-ratio / market-value are
    # all stored as plain text in the source data (even though they look
    # numeric) - the leading "'" forces text entry, same as a human typing
    # it into Excel.
    $q1.Cells.Item($r, 2).Value = "'" + $row[0]
    $q1.Cells.Item($r, 3).Value = $row[1]
    $q1.Cells.Item($r, 4).Value = "'" + $row[2]
    $q1.Cells.Item($r, 5).Value = "'" + $row[3]
    $q1.Cells.Item($r, 6).Value = "'" + $row[4]
    $q1.Cells.Item($r, 7).Value = "'" + $row[5]
    $q1.Cells.Item($r, 8).Value = $row[6]
}

# ---------------------------------------------------------------------------
# Step 2: append a brand-new "总计" sheet right after "2022-Q1" with the
# totals table (old rows shifted down by one, plus the new 2022-Q1 row).
# ---------------------------------------------------------------------------
$total = $wb.Worksheets.Add($null, $q1)
$total.Name = "总计"

$total.Cells.Item(1, 2).Value = "日期"
$total.Cells.Item(1, 3).Value = "持有数量(只)"
$total.Cells.Item(1, 4).Value = "持有市值(亿元)"
$total.Range("B1:D1").Font.Bold = $true
$total.Range("B1:D1").Borders.LineStyle = 1
$total.Range("B1:D1").HorizontalAlignment = -4108
$total.Range("B1:D1").VerticalAlignment = -4160

$totalRows = @(
    ,@("2022-Q1", 12, 4.3)
    ,@("2021-Q4", 13, 4.14)
    ,@("2021-Q3", 7, 2.28)
    ,@("2021-Q2", 12, 2.16)
    ,@("2021-Q1", 5, 0.58)
)

for ($i = 0; $i -lt $totalRows.Length; $i++) {
    $r = $i + 2
    $row = $totalRows[$i]
    $total.Cells.Item($r, 1).Value = $i
    $total.Cells.Item($r, 1).Font.Bold = $true
    $total.Cells.Item($r, 1).Borders.LineStyle = 1
    $total.Cells.Item($r, 1).HorizontalAlignment = -4108
    $total.Cells.Item($r, 1).VerticalAlignment = -4160
    $total.Cells.Item($r, 2).Value = $row[0]
    $total.Cells.Item($r, 3).Value = $row[1]
    $total.Cells.Item($r, 4).Value = $row[2]
}
